# Apply the data edits described by the diff.
$wb = $excel.ActiveWorkbook

# --- Controls sheet ---
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Range("B2").Value = 100          # n_sims: 300 -> 100
$wsControls.Range("B4").Value = 300          # N_1: 10000000 -> 300
$wsControls.Activate()
$wsControls.Range("B3").Select()

# --- Maturity_At_Age sheet ---
$wsMaturity = $wb.Worksheets.Item("Maturity_At_Age")
$wsMaturity.Range("D2").Value = 0.05         # a3 maturity at Sex=1: 0 -> 0.05
$wsMaturity.Activate()
$wsMaturity.Range("D3").Select()

# --- Recruitment_Mortality sheet ---
$wsRecruit = $wb.Worksheets.Item("Recruitment_Mortality")
$wsRecruit.Range("B4").Value = 1             # sigma_rec: 0.75 -> 1
$wsRecruit.Activate()
$wsRecruit.Range("B5").Select()

# Re-select the Controls sheet last, matching tabSelected="1" on that sheet.
$wsControls.Activate()
